# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# table with the latest scraped text values. Every cell on this sheet is
# stored as literal text (not a number), even ones that look numeric
# (e.g. "299.69"), so where the new value would otherwise be auto-typed
# as a number by Excel we write it with a leading apostrophe to force
# text, then reset the cell style back to Normal so no stray
# "quote prefix" formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.942.87'
$ws.Range('E2').Value = '  -0.62%  '
$ws.Range('D3').Value = '2.298.37'
$ws.Range('E3').Value = '  -1.05%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = "'299.69"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.20%  '
$ws.Range('D6').Value = "'97.47"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.30%  '
$ws.Range('E7').Value = '  +1.64%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  -1.87%  '
$ws.Range('D10').Value = "'35.75"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.88%  '
$ws.Range('D11').Value = "'0.0788"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.63%  '
$ws.Range('E12').Value = '  +0.49%  '
$ws.Range('D13').Value = "'17.72"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.07%  '
$ws.Range('D14').Value = "'6.77"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.29%  '
$ws.Range('D15').Value = '2.655.75'
$ws.Range('E15').Value = '  -1.04%  '
$ws.Range('D16').Value = '2.299.89'
$ws.Range('E16').Value = '  +0.32%  '
$ws.Range('D18').Value = '42.890.03'
$ws.Range('E18').Value = '  -0.49%  '
$ws.Range('D19').Value = "'12.60"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.45%  '
$ws.Range('D20').Value = '0.0₃0907'
$ws.Range('E20').Value = '  -0.60%  '
$ws.Range('E21').Value = '  -2.77%  '
$ws.Range('D22').Value = "'68.11"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.10%  '
$ws.Range('D23').Value = "'241.30"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.56%  '
$ws.Range('E24').Value = '  -1.01%  '
$ws.Range('E26').Value = '  -1.46%  '
$ws.Range('E27').Value = '  -0.27%  '
$ws.Range('D28').Value = "'25.13"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.67%  '
$ws.Range('D29').Value = "'166.49"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.08%  '
$ws.Range('D30').Value = "'2.03"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.76%  '
$ws.Range('E31').Value = '  -1.60%  '
$ws.Range('D32').Value = "'32.91"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.18%  '
$ws.Range('E33').Value = '  +0.10%  '
$ws.Range('E34').Value = '  -3.42%  '
$ws.Range('E35').Value = '  -2.09%  '
$ws.Range('E36').Value = '  -0.50%  '
$ws.Range('D37').Value = "'2.39"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.49%  '
$ws.Range('E38').Value = '  -1.82%  '
$ws.Range('E39').Value = '  -2.40%  '
$ws.Range('E40').Value = '  -3.54%  '
$ws.Range('E41').Value = '  -0.78%  '
$ws.Range('E42').Value = '  +0.27%  '
$ws.Range('D43').Value = '2.001.69'
$ws.Range('E43').Value = '  +0.23%  '
$ws.Range('E44').Value = '  -1.22%  '
$ws.Range('E45').Value = '  -3.93%  '
$ws.Range('D46').Value = "'10.13"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.24%  '
$ws.Range('D47').Value = "'17.28"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.70%  '
$ws.Range('D48').Value = "'2.77"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.37%  '
$ws.Range('D49').Value = "'2.91"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.19%  '
$ws.Range('D50').Value = "'53.44"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.12%  '
$ws.Range('D51').Value = '2.521.37'
$ws.Range('E51').Value = '  -1.06%  '
